$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in missing "TAYs" (column F) and "Age_range" (column G) values
# for the Allegheny rows that were left blank.
foreach ($row in 35..39) {
    $ws.Cells.Item($row, 6).Value = "No"
    $ws.Cells.Item($row, 7).Value = "Any age "
}

# Reflect the author's final cursor/viewport position in the sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G47").Select()
